{"js": "// The five \"Empresa:/Cliente:\" dialogue paragraphs collapse into a\n// single descriptive paragraph about \"La Empresa RentaCar\", and the\n// trailing empty paragraph (which only carried the _GoBack bookmark)\n// is removed, its bookmark folded into the new paragraph.\n//\n// We rebuild paragraph 1 from raw WordprocessingML via insertOoxml so\n// we can reproduce the exact run layout, including the <w:proofErr>\n// spell-check markers that bracket the two \"RentaCar\" runs, then\n// delete the now-unneeded trailing paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\nconst newParagraphXml =\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:jc w:val=\"both\"/>' +\n      '<w:rPr><w:lang w:val=\"es-MX\"/></w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr>' +\n      '<w:t xml:space=\"preserve\">La Empresa </w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' +\n      '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr>' +\n      '<w:t>RentaCar</w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' +\n      '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr>' +\n      '<w:t xml:space=\"preserve\"> se dedica a la renta de autos a personas naturales que desean alquilar un autom\\u00f3vil. Los clientes deben suministrar sus datos personales a la empresa y los servicios extra que necesita le proporcionen. A su vez deben informarles sobre los modelos y precios de los veh\\u00edculos que tienen para la renta. El cliente tiene que proporcionar las caracter\\u00edsticas del alquiler que le interesa, y se espera de este que mantenga una buena relaci\\u00f3n con </w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' +\n      '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr>' +\n      '<w:t>RentaCar</w:t>' +\n    '</w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' +\n      '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr>' +\n      '<w:t xml:space=\"preserve\"> que garantice el cumplimiento de lo pactado. Para obtener el auto que necesita el cliente en un momento dado, puede pedirlo prestado a un amigo o familiar, o acudir a la empresa para hacer una reserva.</w:t>' +\n    '</w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>';\n\nconst packagedOoxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nfirstParagraph.insertOoxml(packagedOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-fetch the paragraph collection (the replace above can re-seat\n// paragraph identities) and drop everything after the rebuilt first\n// paragraph.\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\n\nfor (let i = remaining.items.length - 1; i >= 1; i--) {\n  remaining.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Replace the whole body content: the five \"Empresa/Cliente\" dialogue\n# paragraphs collapse into a single descriptive paragraph about\n# \"La Empresa RentaCar\", keeping the trailing _GoBack bookmark.\n# We build the target paragraph as raw WordprocessingML (so we can also\n# emit the <w:proofErr> spell-check markers around \"RentaCar\" exactly as\n# the authored revision has them) and inject it with Range.InsertXML,\n# which parses OOXML and splices it in place of the selected range.\n\n$d = $word.ActiveDocument\n\n$newParagraphXml = '<w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:lang w:val=\"es-MX\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\">La Empresa </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr><w:t>RentaCar</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\"> se dedica a la renta de autos a personas naturales que desean alquilar un autom\u00f3vil. Los clientes deben suministrar sus datos personales a la empresa y los servicios extra que necesita le proporcionen. A su vez deben informarles sobre los modelos y precios de los veh\u00edculos que tienen para la renta. El cliente tiene que proporcionar las caracter\u00edsticas del alquiler que le interesa, y se espera de este que mantenga una buena relaci\u00f3n con </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr><w:t>RentaCar</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:lang w:val=\"es-MX\"/></w:rPr><w:t xml:space=\"preserve\"> que garantice el cumplimiento de lo pactado. Para obtener el auto que necesita el cliente en un momento dado, puede pedirlo prestado a un amigo o familiar, o acudir a la empresa para hacer una reserva.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n\n$wordOpenXml = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Target the full body range (section properties live outside the main\n# range, so this leaves page setup untouched) and splice in the OOXML.\n$d.Range().InsertXML($wordOpenXml)\n"}
